# Refresh the cryptocurrency Price (D) / Volume(1h) (E) snapshot columns
# with the latest scraped values (coinranking.com), cell by cell.
#
# Some Price values are plain decimals (e.g. "0.993"); Excel would
# normally auto-convert a bare numeric-looking string typed into a cell
# into a Number. The source data models these as text, so a leading
# apostrophe (Excel's standard quote-prefix-as-text marker) is used to
# keep them stored as text, matching the original cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.066.70"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").Value = "1.636.89"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("D4").Value = "'0.993"
$ws.Range("E4").Value = "  -0.79%  "
$ws.Range("D5").Value = "'215.45"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("E7").Value = "  -0.86%  "
$ws.Range("E8").Value = "  -1.09%  "
$ws.Range("E9").Value = "  -0.94%  "
$ws.Range("D10").Value = "'19.80"
$ws.Range("E10").Value = "  +0.58%  "
$ws.Range("D11").Value = "'0.0788"
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").Value = "'4.25"
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("D13").Value = "1.862.01"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").Value = "1.633.16"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("E15").Value = "  -1.48%  "
$ws.Range("D16").Value = "0.0₃0763"
$ws.Range("E16").Value = "  -0.71%  "
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").Value = "26.033.90"
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("D19").Value = "'0.993"
$ws.Range("E19").Value = "  -0.87%  "
$ws.Range("D20").Value = "'4.45"
$ws.Range("E20").Value = "  -0.52%  "
$ws.Range("D21").Value = "'192.92"
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").Value = "'6.37"
$ws.Range("E23").Value = "  +1.23%  "
$ws.Range("D24").Value = "'0.992"
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("E25").Value = "  -1.79%  "
$ws.Range("D26").Value = "'141.56"
$ws.Range("E26").Value = "  -0.90%  "
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("D28").Value = "'6.89"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("D29").Value = "'15.62"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  +0.54%  "
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("D33").Value = "'3.24"
$ws.Range("E33").Value = "  -0.66%  "
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("D36").Value = "'0.907"
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").Value = "1.140.51"
$ws.Range("E37").Value = "  +0.92%  "
$ws.Range("D38").Value = "'0.545"
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("E39").Value = "  -1.60%  "
$ws.Range("E40").Value = "  -0.41%  "
$ws.Range("D41").Value = "'0.992"
$ws.Range("E41").Value = "  -0.75%  "
$ws.Range("D42").Value = "'5.56"
$ws.Range("E42").Value = "  -0.61%  "
$ws.Range("D43").Value = "'100.21"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").Value = "'0.796"
$ws.Range("E44").Value = "  -1.38%  "
$ws.Range("D45").Value = "1.772.38"
$ws.Range("D46").Value = "0.0₆0106"
$ws.Range("E46").Value = "  -2.22%  "
$ws.Range("D47").Value = "'55.63"
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("D48").Value = "'0.0515"
$ws.Range("E48").Value = "  +2.36%  "
$ws.Range("E49").Value = "  +4.32%  "
$ws.Range("D50").Value = "'0.415"
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("D51").Value = "'7.63"
$ws.Range("E51").Value = "  +0.98%  "
